$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.124.20'
$ws.Range("E2").Value = '  +6.34%  '
$ws.Range("D3").Value = '3.109.08'
$ws.Range("E3").Value = '  +4.29%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").Value = "'583.67"
$ws.Range("E5").Value = '  +3.89%  '
$ws.Range("D6").Value = "'144.74"
$ws.Range("E6").Value = '  +5.04%  '
$ws.Range("D8").Value = '3.098.72'
$ws.Range("E8").Value = '  +4.25%  '
$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = '  +0.93%  '
$ws.Range("E10").Value = '  +13.74%  '
$ws.Range("E11").Value = '  +7.10%  '
$ws.Range("E12").Value = '  +3.40%  '
$ws.Range("D13").Value = "'0.0000247"
$ws.Range("E13").Value = '  +7.88%  '
$ws.Range("D14").Value = "'35.44"
$ws.Range("E14").Value = '  +5.25%  '
$ws.Range("E15").Value = '  +0.68%  '
$ws.Range("D16").Value = '3.622.83'
$ws.Range("E16").Value = '  +4.39%  '
$ws.Range("E17").Value = '  +1.72%  '
$ws.Range("D18").Value = '63.064.49'
$ws.Range("E18").Value = '  +6.32%  '
$ws.Range("D19").Value = '3.107.64'
$ws.Range("E19").Value = '  +4.44%  '
$ws.Range("D20").Value = "'465.05"
$ws.Range("E20").Value = '  +6.37%  '
$ws.Range("D21").Value = "'14.18"
$ws.Range("E21").Value = '  +4.72%  '
$ws.Range("E22").Value = '  +1.16%  '
$ws.Range("D23").Value = "'7.52"
$ws.Range("E23").Value = '  +7.12%  '
$ws.Range("D24").Value = "'13.28"
$ws.Range("E24").Value = '  -0.42%  '
$ws.Range("D25").Value = "'82.03"
$ws.Range("E25").Value = '  +2.68%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").Value = "'8.49"
$ws.Range("E27").Value = '  +9.75%  '
$ws.Range("E28").Value = '  +0.26%  '
$ws.Range("B29").Value = 'FirstDigitalUSD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = '  +0.29%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = "'2.67"
$ws.Range("E30").Value = '  +5.23%  '
$ws.Range("D31").Value = "'6.85"
$ws.Range("E31").Value = '  +9.86%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = "'26.92"
$ws.Range("E32").Value = '  +4.81%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = "'0.111"
$ws.Range("E33").Value = '  +4.68%  '
$ws.Range("D34").Value = '0.0₃0868'
$ws.Range("E34").Value = '  +13.53%  '
$ws.Range("E35").Value = '  +15.82%  '
$ws.Range("D36").Value = "'1.05"
$ws.Range("E36").Value = '  +6.61%  '
$ws.Range("D37").Value = "'6.06"
$ws.Range("E37").Value = '  +2.95%  '
$ws.Range("D38").Value = "'3.30"
$ws.Range("E38").Value = '  +18.84%  '
$ws.Range("E39").Value = '  +4.62%  '
$ws.Range("D40").Value = "'438.61"
$ws.Range("E40").Value = '  +9.74%  '
$ws.Range("D41").Value = "'8.72"
$ws.Range("E41").Value = '  +0.61%  '
$ws.Range("D42").Value = '2.914.37'
$ws.Range("E42").Value = '  +6.09%  '
$ws.Range("D43").Value = "'0.0368"
$ws.Range("E43").Value = '  +4.78%  '
$ws.Range("E44").Value = '  +12.30%  '
$ws.Range("D45").Value = "'0.111"
$ws.Range("E45").Value = '  +6.71%  '
$ws.Range("E46").Value = '  +7.83%  '
$ws.Range("D47").Value = "'35.10"
$ws.Range("E47").Value = '  +0.59%  '
$ws.Range("D49").Value = "'122.81"
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("E50").Value = '  +0.52%  '
$ws.Range("D51").Value = "'24.52"
$ws.Range("E51").Value = '  +5.15%  '

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
